$d = $word.ActiveDocument
$t = $d.Tables(1)

# Remove the row that is deleted in this revision (original row 7:
#   27+1=, 88-23=, 95+3=, 38+53=, 79+9=)
$t.Rows(7).Delete()

# Update all remaining problem cells to their new values
$t.Cell(1, 1).Range.Text = "36-22="
$t.Cell(1, 2).Range.Text = "26-7="
$t.Cell(1, 3).Range.Text = "54-47="
$t.Cell(1, 4).Range.Text = "10+12="
$t.Cell(1, 5).Range.Text = "22+2="
$t.Cell(2, 1).Range.Text = "98-36="
$t.Cell(2, 2).Range.Text = "83-48="
$t.Cell(2, 3).Range.Text = "26+41="
$t.Cell(2, 4).Range.Text = "30+68="
$t.Cell(2, 5).Range.Text = "25+42="
$t.Cell(3, 1).Range.Text = "46+17="
$t.Cell(3, 2).Range.Text = "32-15="
$t.Cell(3, 3).Range.Text = "21+67="
$t.Cell(3, 4).Range.Text = "8+45="
$t.Cell(3, 5).Range.Text = "34+19="
$t.Cell(4, 1).Range.Text = "48+8="
$t.Cell(4, 2).Range.Text = "4+35="
$t.Cell(4, 3).Range.Text = "31+21="
$t.Cell(4, 4).Range.Text = "75-62="
$t.Cell(4, 5).Range.Text = "93-83="
$t.Cell(5, 1).Range.Text = "23-20="
$t.Cell(5, 2).Range.Text = "56-0="
$t.Cell(5, 3).Range.Text = "2+56="
$t.Cell(5, 4).Range.Text = "4+33="
$t.Cell(5, 5).Range.Text = "81-59="
$t.Cell(6, 1).Range.Text = "75-20="
$t.Cell(6, 2).Range.Text = "2+28="
$t.Cell(6, 3).Range.Text = "27-24="
$t.Cell(6, 4).Range.Text = "23-0="
$t.Cell(6, 5).Range.Text = "49-23="
$t.Cell(7, 1).Range.Text = "24+18="
$t.Cell(7, 3).Range.Text = "1+17="
$t.Cell(7, 4).Range.Text = "97-83="
$t.Cell(7, 5).Range.Text = "45+12="
$t.Cell(8, 1).Range.Text = "6+47="
$t.Cell(8, 2).Range.Text = "81-56="
$t.Cell(8, 3).Range.Text = "53-15="
$t.Cell(8, 4).Range.Text = "44+15="
$t.Cell(8, 5).Range.Text = "19+17="
$t.Cell(9, 1).Range.Text = "40+56="
$t.Cell(9, 2).Range.Text = "60+31="
$t.Cell(9, 3).Range.Text = "27+23="
$t.Cell(9, 4).Range.Text = "63-6="
$t.Cell(9, 5).Range.Text = "25+42="
$t.Cell(10, 1).Range.Text = "30+9="
$t.Cell(10, 2).Range.Text = "46+7="
$t.Cell(10, 3).Range.Text = "9+23="
$t.Cell(10, 4).Range.Text = "21+36="
$t.Cell(10, 5).Range.Text = "1+54="
$t.Cell(11, 1).Range.Text = "18+4="
$t.Cell(11, 2).Range.Text = "26-26="
$t.Cell(11, 3).Range.Text = "29+0="
$t.Cell(11, 4).Range.Text = "22+7="
$t.Cell(11, 5).Range.Text = "78-77="
$t.Cell(12, 1).Range.Text = "0+97="
$t.Cell(12, 2).Range.Text = "66-33="
$t.Cell(12, 3).Range.Text = "17+57="
$t.Cell(12, 4).Range.Text = "40+46="
$t.Cell(12, 5).Range.Text = "17+30="
$t.Cell(13, 1).Range.Text = "16-6="
$t.Cell(13, 2).Range.Text = "63+5="
$t.Cell(13, 3).Range.Text = "2+77="
$t.Cell(13, 4).Range.Text = "17-4="
$t.Cell(13, 5).Range.Text = "64-23="
$t.Cell(14, 1).Range.Text = "80-56="
$t.Cell(14, 2).Range.Text = "67-43="
$t.Cell(14, 3).Range.Text = "60-20="
$t.Cell(14, 4).Range.Text = "16-8="
$t.Cell(14, 5).Range.Text = "18+34="
$t.Cell(15, 1).Range.Text = "30-16="
$t.Cell(15, 2).Range.Text = "8+32="
$t.Cell(15, 3).Range.Text = "66+5="
$t.Cell(15, 4).Range.Text = "85+7="
$t.Cell(15, 5).Range.Text = "27+24="
$t.Cell(16, 1).Range.Text = "17+18="
$t.Cell(16, 2).Range.Text = "36+43="
$t.Cell(16, 3).Range.Text = "91-37="
$t.Cell(16, 4).Range.Text = "58-36="
$t.Cell(16, 5).Range.Text = "50+39="
$t.Cell(17, 1).Range.Text = "48-46="
$t.Cell(17, 2).Range.Text = "39-2="
$t.Cell(17, 3).Range.Text = "45-28="
$t.Cell(17, 4).Range.Text = "84-79="
$t.Cell(17, 5).Range.Text = "89-50="
$t.Cell(18, 1).Range.Text = "48-32="
$t.Cell(18, 2).Range.Text = "11+52="
$t.Cell(18, 3).Range.Text = "48+43="
$t.Cell(18, 4).Range.Text = "49-21="
$t.Cell(18, 5).Range.Text = "54-33="
$t.Cell(19, 1).Range.Text = "81-57="
$t.Cell(19, 2).Range.Text = "13-12="
$t.Cell(19, 3).Range.Text = "44-21="
$t.Cell(19, 4).Range.Text = "56+41="
$t.Cell(19, 5).Range.Text = "63+3="

# Append a brand-new row at the end of the table
$newRow = $t.Rows.Add()
$newRow.Cells(1).Range.Text = "71-20="
$newRow.Cells(2).Range.Text = "18+78="
$newRow.Cells(3).Range.Text = "72-4="
$newRow.Cells(4).Range.Text = "43+35="
$newRow.Cells(5).Range.Text = "14+28="
